$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add the new mail-log entry as row 16
$ws.Range("A16").Value = "Offerte voor zakelijke samenwerking"
$ws.Range("B16").Value = "mailmind.test@zohomail.eu"
$ws.Range("C16").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$ws.Range("D16").Value = "Offerte / Prijsaanvraag"
$ws.Range("F16").Value = "2025-06-20 15:00:10"
$ws.Range("G16").Value = "Nee"

# Extend the conditional formatting ranges on columns D and G to include the new row
$dFormats = $ws.Range("D2:D15").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($ws.Range("D2:D16"))
}

$gFormats = $ws.Range("G2:G15").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($ws.Range("G2:G16"))
}

# Update the Dashboard count for "Offerte / Prijsaanvraag" from 2 to 3
$dash.Range("B4").Value = 3
